$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other header cells (copy format from G1).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Add the data value for the new column (plain, unstyled like the other data cells)
$ws.Range("H2").Value = 0
